# StructureDefinition-eclaire-review-date.xlsx
# "changement type date to instant pour extension ECLAIREReviewDate (#50)"
#
# 1. Bump the resource's "Date" metadata field (Metadata!B8) to the new
#    last-modified timestamp.
# 2. Change the data type of the Extension.value[x] element (Elements!K5)
#    from "date" to "instant".
# 3. The FHIR IG Publisher export regenerates the root Extension row's
#    Invariants column (Elements!AJ1) once any child element changes -
#    it picks up the ele-1 / ext-1 invariants that already apply to the
#    Extension type.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-08-23T14:17:04+00:00"

$elements = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type column: date -> instant
$elements.Range("K5").Value = "instant`n"

# Root Extension row Invariants column gets populated with the
# ele-1 / ext-1 constraint text (already shown on the Extension.extension row).
$elements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
